$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misparsed as a number/date by Excel's
# input parser (e.g. "1.000", "242.45") are forced to Text format first
# so the literal string is preserved exactly, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.333.77'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.38'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7103'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.45'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08018'
$ws.Range("E8").Value = '  +3.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3140'
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08332'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.881.95'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.276'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.79'
$ws.Range("E14").Value = '  +3.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7175'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.377'
$ws.Range("E16").Value = '  +5.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008755'
$ws.Range("E17").Value = '  +6.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.360.65'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.60'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.29'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.138.82'
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.851'
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1573'
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.58'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.068'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.60'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.428'
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.344'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.206'
$ws.Range("E32").Value = '  -5.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05415'
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.940'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7762'
$ws.Range("E35").Value = '  +3.55%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01886'
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.276.08'
$ws.Range("E39").Value = '  +5.86%  '
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.548'
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9206'
$ws.Range("E42").Value = '  +3.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '112.74'
$ws.Range("E43").Value = '  +4.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '74.57'
$ws.Range("E44").Value = '  +2.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.031.05'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000127'
$ws.Range("E47").Value = '  +3.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.807'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5218'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.555'
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4375'
$ws.Range("E51").Value = '  +1.27%  '
